# Project creation - caso de prueba de fecha de inicio
#
# Adds two new worksheets (Hoja11, Hoja12) that exercise the "project
# start date" test case, re-using the existing Username/Contraseña
# ("pepeusername" / "P4ssword.") test credentials from Hoja10, and
# updates the selection state on Hoja10 that results from the tab
# switch.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add "Hoja11" right after the last existing sheet (Hoja10).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$hoja11 = $wb.Worksheets.Add($null, $lastSheet)
$hoja11.Name = "Hoja11"

$hoja11.Range("A1").Value = "Username"
$hoja11.Range("B1").Value = "Contraseña"
$hoja11.Range("C1").Value = "nombre proyecto"
$hoja11.Range("D1").Value = "Descripcion"
$hoja11.Range("E1").Value = "Fecha de inicio"

$hoja11.Range("A2").Value = "pepeusername"
$hoja11.Range("B2").Value = "P4ssword."
$hoja11.Range("C2").Value = "Date Test 1"
$hoja11.Range("D2").Value = "project created in order to test the start date"

$hoja11.Range("A3").Value = "pepeusername"
$hoja11.Range("B3").Value = "P4ssword."
$hoja11.Range("C3").Value = "Date Test12"
$hoja11.Range("D3").Value = "project created in order to test the start date"
$hoja11.Range("E3").Value = "13 October 2021"

[void]$hoja11.Range("A1:E1").Select()

# ---------------------------------------------------------------------
# 2. Add "Hoja12" right after "Hoja11".
# ---------------------------------------------------------------------
$hoja12 = $wb.Worksheets.Add($null, $hoja11)
$hoja12.Name = "Hoja12"

$hoja12.Range("A1").Value = "Username"
$hoja12.Range("B1").Value = "Contraseña"
$hoja12.Range("C1").Value = "nombre proyecto"
$hoja12.Range("D1").Value = "Descripcion"
$hoja12.Range("E1").Value = "Fecha de inicio"

$hoja12.Range("A2").Value = "pepeusername"
$hoja12.Range("B2").Value = "P4ssword."
$hoja12.Range("C2").Value = "Date Test12"
$hoja12.Range("D2").Value = "project created in order to test the start date"
$hoja12.Range("E2").Value = "13 October 2021"

# Hoja12 ends up being the active sheet/tab (last one added), matching
# activeTab and the lack of an explicit tabSelected elsewhere.

# ---------------------------------------------------------------------
# 3. Hoja10 loses the tab-selected state (Hoja12 is now active) and its
#    selection becomes a multi-cell block.
# ---------------------------------------------------------------------
$hoja10 = $wb.Worksheets.Item("Hoja10")
[void]$hoja10.Range("A1:E3").Select()

# Re-activate Hoja12 so it is the workbook's active sheet/tab.
$hoja12.Activate()
